# cicada_ig/temp/pages/StructureDefinition-inadvertent-administration-status.xlsx
#
# This IG-generated page was regenerated for the "cicada" IG (was "pythia"),
# with a new generation timestamp, and a new "Jurisdiction" metadata row.
#
# Changes:
#   1. Metadata!B2  - canonical URL: pythia -> cicada
#   2. Metadata!B8  - generation Date stamp updated
#   3. Metadata     - new "Jurisdiction" row inserted after "Contact" (row 10),
#                      pushing Description/Purpose/.../Context down by one row
#                      (dimension grows from B20 to B21)
#   4. Elements!R5  - "Fixed Value" for Extension.url mirrors the page's own
#                      canonical URL (same underlying string as Metadata!B2),
#                      so it needs the same pythia -> cicada update

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# 1. Canonical URL
$meta.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/inadvertent-administration-status"

# 2. Generation date
$meta.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10) and
#    before "Description" (row 11), carrying the same row styling.
$meta.Rows.Item(11).Insert()
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# 4. Keep the Elements grid's mirrored "Fixed Value" for Extension.url in sync
#    with the page's own canonical URL.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("R5").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/inadvertent-administration-status"
